$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Table2" question columns I:M, mirroring the pattern already used
# for the Table (A,B,C,D) columns F:H, but with a 3-row repeating cycle.
$header = "Table2"

$col1Values = @("1,3,5;OO;1", "1,2,4;OO;2", "1,4;OO;3")
$col2Values = @("1,2,3;OO;2", "1;OO;3", "1,2;OO;4")

$columns = @(9, 10, 11, 12, 13)  # I, J, K, L, M

foreach ($col in $columns) {
    # Header
    $ws.Cells.Item(1, $col).Value = $header

    # choose the value cycle depending on whether this is an "odd" (I,K,M)
    # or "even" (J,L) column in the pattern
    if ((($col - 9) % 2) -eq 0) {
        $values = $col1Values
    } else {
        $values = $col2Values
    }

    for ($r = 2; $r -le 14; $r++) {
        $idx = ($r - 2) % 3
        $ws.Cells.Item($r, $col).Value = $values[$idx]
    }
}
